# Update res_line/pl_mw.xlsx data table (rows 2-25, columns B,C,E,F,G,H,J,K,L,N,O)
# for the "case with 380 kV" re-run. Columns A, D, I, M are untouched (structural
# zero columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6313203717040778
$ws.Range("C2").Value = 0.298564395615148
$ws.Range("E2").Value = 0.2776979355111315
$ws.Range("F2").Value = 1.886419488441902
$ws.Range("G2").Value = 0.6396534143365358
$ws.Range("H2").Value = 0.7680410200398313
$ws.Range("J2").Value = 0.04602985507985835
$ws.Range("K2").Value = 0.2622393503981471
$ws.Range("L2").Value = 0.4667620490305353
$ws.Range("N2").Value = 1.657807643039227
$ws.Range("O2").Value = 2.804474144802967
$ws.Range("B3").Value = 0.590842794766985
$ws.Range("C3").Value = 0.3005191531905909
$ws.Range("E3").Value = 0.2764103369528108
$ws.Range("F3").Value = 1.885398137126714
$ws.Range("G3").Value = 0.643612652319618
$ws.Range("H3").Value = 0.7732198337806011
$ws.Range("J3").Value = 0.04455423301908468
$ws.Range("K3").Value = 0.2310627759570707
$ws.Range("L3").Value = 0.4567322239196159
$ws.Range("N3").Value = 1.672942198728789
$ws.Range("O3").Value = 2.82348138757159
$ws.Range("B4").Value = 0.5661547638283082
$ws.Range("C4").Value = 0.3017852055618193
$ws.Range("E4").Value = 0.2757278429815209
$ws.Range("F4").Value = 1.885664671379807
$ws.Range("G4").Value = 0.6463966460535389
$ws.Range("H4").Value = 0.7766748148483202
$ws.Range("J4").Value = 0.04364000301762516
$ws.Range("K4").Value = 0.2118920633115664
$ws.Range("L4").Value = 0.4507589251714137
$ws.Range("N4").Value = 1.6827377900551
$ws.Range("O4").Value = 2.836468379794894
$ws.Range("B5").Value = 0.5561365717192643
$ws.Range("C5").Value = 0.3023177109638331
$ws.Range("E5").Value = 0.2754769834422461
$ws.Range("F5").Value = 1.885998330104819
$ws.Range("G5").Value = 0.64761992367891
$ws.Range("H5").Value = 0.7781520178016166
$ws.Range("J5").Value = 0.04326540162008641
$ws.Range("K5").Value = 0.2040732881050076
$ws.Range("L5").Value = 0.4483715352459114
$ws.Range("N5").Value = 1.686856095453177
$ws.Range("O5").Value = 2.842091900769333
$ws.Range("B6").Value = 0.5544756417487804
$ws.Range("C6").Value = 0.3024071351915172
$ws.Range("E6").Value = 0.2754369773474536
$ws.Range("F6").Value = 1.886067335638792
$ws.Range("G6").Value = 0.6478284102833385
$ws.Range("H6").Value = 0.7784014922319997
$ws.Range("J6").Value = 0.04320307624668018
$ws.Range("K6").Value = 0.2027746061643967
$ws.Range("L6").Value = 0.4479779434142728
$ws.Range("N6").Value = 1.687547580395217
$ws.Range("O6").Value = 2.843045692079329
$ws.Range("B7").Value = 0.5660194823288691
$ws.Range("C7").Value = 0.3017923199526251
$ws.Range("E7").Value = 0.2757243493109627
$ws.Range("F7").Value = 1.88566825958749
$ws.Range("G7").Value = 0.6464127841311438
$ws.Range("H7").Value = 0.7766944563505263
$ws.Range("J7").Value = 0.04363495927362493
$ws.Range("K7").Value = 0.2117866424107291
$ws.Range("L7").Value = 0.4507265382701178
$ws.Range("N7").Value = 1.682792818595328
$ws.Range("O7").Value = 2.836542879245272
$ws.Range("B8").Value = 0.6173298475558795
$ws.Range("C8").Value = 0.2992247507552204
$ws.Range("E8").Value = 0.2772315834021875
$ws.Range("F8").Value = 1.88588204235154
$ws.Range("G8").Value = 0.6409452995607126
$ws.Range("H8").Value = 0.7697696098423421
$ws.Range("J8").Value = 0.04552276933087285
$ws.Range("K8").Value = 0.251495822350023
$ws.Range("L8").Value = 0.463265500010877
$ws.Range("N8").Value = 1.662921695922943
$ws.Range("O8").Value = 2.81075471555063
$ws.Range("B9").Value = 0.7192315344318843
$ws.Range("C9").Value = 0.2947107695615294
$ws.Range("E9").Value = 0.2810418237949861
$ws.Range("F9").Value = 1.89337970033138
$ws.Range("G9").Value = 0.6330243802759483
$ws.Range("H9").Value = 0.7583700014021062
$ws.Range("J9").Value = 0.04915926357849898
$ws.Range("K9").Value = 0.3291222974670802
$ws.Range("L9").Value = 0.4893135578579404
$ws.Range("N9").Value = 1.62794166558766
$ws.Range("O9").Value = 2.770622740548532
$ws.Range("B10").Value = 0.7948478832560113
$ws.Range("C10").Value = 0.2917101544301692
$ws.Range("E10").Value = 0.2843585320033384
$ws.Range("F10").Value = 1.903190160455793
$ws.Range("G10").Value = 0.6289126755423098
$ws.Range("H10").Value = 0.7513192384409848
$ws.Range("J10").Value = 0.05179062657331457
$ws.Range("K10").Value = 0.3859853364062928
$ws.Range("L10").Value = 0.5093307200364734
$ws.Range("N10").Value = 1.604667426207982
$ws.Range("O10").Value = 2.747493155325188
$ws.Range("B11").Value = 0.8294039500244423
$ws.Range("C11").Value = 0.2904132837395954
$ws.Range("E11").Value = 0.2859789839355003
$ws.Range("F11").Value = 1.908584701588723
$ws.Range("G11").Value = 0.6274131013666491
$ws.Range("H11").Value = 0.7483983444640216
$ws.Range("J11").Value = 0.0529788526689785
$ws.Range("K11").Value = 0.4118129564347157
$ws.Range("L11").Value = 0.5186260866237689
$ws.Range("N11").Value = 1.594604993209554
$ws.Range("O11").Value = 2.738349292313217
$ws.Range("B12").Value = 0.8425114099506743
$ws.Range("C12").Value = 0.2899319610645783
$ws.Range("E12").Value = 0.2866085904902249
$ws.Range("F12").Value = 1.910761171573057
$ws.Range("G12").Value = 0.6268985832736718
$ws.Range("H12").Value = 0.7473334088479078
$ws.Range("J12").Value = 0.05342752552846974
$ws.Range("K12").Value = 0.4215870292771342
$ws.Range("L12").Value = 0.5221730256412656
$ws.Range("N12").Value = 1.590870071729842
$ws.Range("O12").Value = 2.735084749000293
$ws.Range("B13").Value = 0.839687526939116
$ws.Range("C13").Value = 0.290035188157165
$ws.Range("E13").Value = 0.2864722841773428
$ws.Range("F13").Value = 1.910286489490119
$ws.Range("G13").Value = 0.6270070215901171
$ws.Range("H13").Value = 0.7475609328624415
$ws.Range("J13").Value = 0.05333095313152825
$ws.Range("K13").Value = 0.4194822966795186
$ws.Range("L13").Value = 0.521407932801182
$ws.Range("N13").Value = 1.591671095995576
$ws.Range("O13").Value = 2.735779020585625
$ws.Range("B14").Value = 0.8304818761468198
$ws.Range("C14").Value = 0.2903734892699514
$ws.Range("E14").Value = 0.2860304623139314
$ws.Range("F14").Value = 1.90876108380337
$ws.Range("G14").Value = 0.6273697026297782
$ws.Range("H14").Value = 0.7483099073319863
$ws.Range("J14").Value = 0.05301579106322407
$ws.Range("K14").Value = 0.4126172046170211
$ws.Range("L14").Value = 0.5189173560866607
$ws.Range("N14").Value = 1.594296205819028
$ws.Range("O14").Value = 2.738076747951794
$ws.Range("B15").Value = 0.8248459661049878
$ws.Range("C15").Value = 0.2905819806969241
$ws.Range("E15").Value = 0.2857619118900416
$ws.Range("F15").Value = 1.90784412804355
$ws.Range("G15").Value = 0.6275988016861476
$ws.Range("H15").Value = 0.7487740321617693
$ws.Range("J15").Value = 0.05282257763084175
$ws.Range("K15").Value = 0.4084113015341018
$ws.Range("L15").Value = 0.5173953125489845
$ws.Range("N15").Value = 1.595913993907374
$ws.Range("O15").Value = 2.739509959841072
$ws.Range("B16").Value = 0.7925926645926324
$ws.Range("C16").Value = 0.2917962769208717
$ws.Range("E16").Value = 0.2842548717455458
$ws.Range("F16").Value = 1.902856336931762
$ws.Range("G16").Value = 0.6290181396700589
$ws.Range("H16").Value = 0.7515158864126903
$ws.Range("J16").Value = 0.05171279470893353
$ws.Range("K16").Value = 0.3842965927626949
$ws.Range("L16").Value = 0.5087270362835596
$ws.Range("N16").Value = 1.605335595995463
$ws.Range("O16").Value = 2.748118446514752
$ws.Range("B17").Value = 0.7728461121099315
$ws.Range("C17").Value = 0.2925586407884762
$ws.Range("E17").Value = 0.2833588940916485
$ws.Range("F17").Value = 1.900034906335961
$ws.Range("G17").Value = 0.629983850302537
$ws.Range("H17").Value = 0.753271270741088
$ws.Range("J17").Value = 0.05102971369637999
$ws.Range("K17").Value = 0.3694924268644115
$ws.Range("L17").Value = 0.5034576708196568
$ws.Range("N17").Value = 1.611249935324924
$ws.Range("O17").Value = 2.753752315789527
$ws.Range("B18").Value = 0.7615033230237316
$ws.Range("C18").Value = 0.2930035457948303
$ws.Range("E18").Value = 0.2828540666894241
$ws.Range("F18").Value = 1.898499809700496
$ws.Range("G18").Value = 0.6305742082067454
$ws.Range("H18").Value = 0.7543078948035316
$ws.Range("J18").Value = 0.05063599690537757
$ws.Range("K18").Value = 0.3609737639060597
$ws.Range("L18").Value = 0.5004447194628767
$ws.Range("N18").Value = 1.61470113554633
$ws.Range("O18").Value = 2.757122469049534
$ws.Range("B19").Value = 0.7576654359591828
$ws.Range("C19").Value = 0.2931552851830972
$ws.Range("E19").Value = 0.2826849492645209
$ws.Range("F19").Value = 1.897995126867073
$ws.Range("G19").Value = 0.6307800882336352
$ws.Range("H19").Value = 0.7546635124349592
$ws.Range("J19").Value = 0.05050254973120261
$ws.Range("K19").Value = 0.3580888735814654
$ws.Range("L19").Value = 0.4994276592311735
$ws.Range("N19").Value = 1.615878142418573
$ws.Range("O19").Value = 2.758285824338316
$ws.Range("B20").Value = 0.7749466299238748
$ws.Range("C20").Value = 0.2924768222438718
$ws.Range("E20").Value = 0.2834531847381498
$ws.Range("F20").Value = 1.900326176406878
$ws.Range("G20").Value = 0.6298774359807098
$ws.Range("H20").Value = 0.7530816158674085
$ws.Range("J20").Value = 0.05110251452224901
$ws.Range("K20").Value = 0.3710687425792685
$ws.Range("L20").Value = 0.504016758409449
$ws.Range("N20").Value = 1.610615228207127
$ws.Range("O20").Value = 2.753139157988159
$ws.Range("B21").Value = 0.8331852134424196
$ws.Range("C21").Value = 0.290273856940904
$ws.Range("E21").Value = 0.2861598031491894
$ws.Range("F21").Value = 1.909205506984819
$ws.Range("G21").Value = 0.6272617267686016
$ws.Range("H21").Value = 0.7480887992246466
$ws.Range("J21").Value = 0.05310839675525614
$ws.Range("K21").Value = 0.4146338238234364
$ws.Range("L21").Value = 0.5196481685273682
$ws.Range("N21").Value = 1.593523098356322
$ws.Range("O21").Value = 2.737396475809646
$ws.Range("B22").Value = 0.8713743192980701
$ws.Range("C22").Value = 0.2888910526964752
$ws.Range("E22").Value = 0.2880218209114318
$ws.Range("F22").Value = 1.915787676560171
$ws.Range("G22").Value = 0.625863108603852
$ws.Range("H22").Value = 0.7450654935597356
$ws.Range("J22").Value = 0.05441187246029244
$ws.Range("K22").Value = 0.4430692107370646
$ws.Range("L22").Value = 0.5300213699196377
$ws.Range("N22").Value = 1.582792492669604
$ws.Range("O22").Value = 2.728261999518764
$ws.Range("B23").Value = 0.850980756390328
$ws.Range("C23").Value = 0.2896238766456829
$ws.Range("E23").Value = 0.2870195362934211
$ws.Range("F23").Value = 1.912203470113951
$ws.Range("G23").Value = 0.6265811273874675
$ws.Range("H23").Value = 0.7466571676208389
$ws.Range("J23").Value = 0.05371687374413625
$ws.Range("K23").Value = 0.4278962796825283
$ws.Range("L23").Value = 0.5244707015349377
$ws.Range("N23").Value = 1.588479353439929
$ws.Range("O23").Value = 2.733031658603977
$ws.Range("B24").Value = 0.7739969553310004
$ws.Range("C24").Value = 0.2925137918109408
$ws.Range("E24").Value = 0.2834105239035765
$ws.Range("F24").Value = 1.900194222246512
$ws.Range("G24").Value = 0.62992543638196
$ws.Range("H24").Value = 0.7531672733878381
$ws.Range("J24").Value = 0.05106960439641739
$ws.Range("K24").Value = 0.3703561137205043
$ws.Range("L24").Value = 0.5037639435674635
$ws.Range("N24").Value = 1.610902020532876
$ws.Range("O24").Value = 2.753415958047739
$ws.Range("B25").Value = 0.6915305162949323
$ws.Range("C25").Value = 0.2958763316969626
$ws.Range("E25").Value = 0.279919941963378
$ws.Range("F25").Value = 1.890594958076107
$ws.Range("G25").Value = 0.6348673164085739
$ws.Range("H25").Value = 0.7612209761584907
$ws.Range("J25").Value = 0.04818255231369406
$ws.Range("K25").Value = 0.3081504693129489
$ws.Range("L25").Value = 0.5173953125489845
$ws.Range("N25").Value = 1.636978373174564
$ws.Range("O25").Value = 2.780362838313096
